$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.881.92"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.630.34"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.52"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "23.27"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("E10").Value = "  -0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0880"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.862.61"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.627.98"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.01"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -0.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.556"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("E16").Value = "  -1.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.902.60"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "228.02"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("E19").Value = "  -0.18%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0718"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -0.69%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("E22").Value = "  -0.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.90"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  -4.70%  "
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "155.34"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.91"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  -0.31%  "
$ws.Range("B28").Value = "BinanceUSD"
$ws.Range("C28").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.08%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.42"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  -0.98%  "
$ws.Range("E30").Value = "  -0.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0480"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.24%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.416.75"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  +1.22%  "
$ws.Range("E35").Value = "  +2.88%  "
$ws.Range("E36").Value = "  -3.92%  "
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("E38").Value = "  -0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.552"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.01"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -1.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "65.87"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("E43").Value = "  -0.72%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.41"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.58%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.771.69"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.11"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -3.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "88.59"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.62"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +1.50%  "
$ws.Range("E51").Value = "  -0.11%  "
